$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill the full data table (header + 10 data rows)
$ws.Range("A1").Value = "Testcase name"
$ws.Range("B1").Value = "Insurancestatus"
$ws.Range("C1").Value = "car"
$ws.Range("D1").Value = "model"
$ws.Range("E1").Value = "type"
$ws.Range("F1").Value = "carhsn"
$ws.Range("A2").Value = "BMW 1er Model"
$ws.Range("B2").Value = "Modify"
$ws.Range("C2").Value = "BMW"
$ws.Range("D2").Value = "1er"
$ws.Range("E2").Value = "Cabrio"
$ws.Range("F2").Value = "null"
$ws.Range("A3").Value = "BMW X1 Model"
$ws.Range("B3").Value = "Purchase"
$ws.Range("C3").Value = "BMW"
$ws.Range("D3").Value = "X1"
$ws.Range("E3").Value = "null"
$ws.Range("F3").Value = "null"
$ws.Range("A4").Value = "BMW X3 Model"
$ws.Range("B4").Value = "Purchase"
$ws.Range("C4").Value = "BMW"
$ws.Range("D4").Value = "X3"
$ws.Range("E4").Value = "null"
$ws.Range("F4").Value = "null"
$ws.Range("A5").Value = "AUDI A1  Model"
$ws.Range("B5").Value = "Purchase"
$ws.Range("C5").Value = "AUDI"
$ws.Range("D5").Value = "A1"
$ws.Range("E5").Value = "null"
$ws.Range("F5").Value = "null"
$ws.Range("A6").Value = "AUDI A4  Model"
$ws.Range("B6").Value = "Purchase"
$ws.Range("C6").Value = "AUDI"
$ws.Range("D6").Value = "A4"
$ws.Range("E6").Value = "Kombi"
$ws.Range("F6").Value = "null"
$ws.Range("A7").Value = "AUDI A5  Model"
$ws.Range("B7").Value = "Purchase"
$ws.Range("C7").Value = "AUDI"
$ws.Range("D7").Value = "A5"
$ws.Range("E7").Value = "Cabrio"
$ws.Range("F7").Value = "null"
$ws.Range("A8").Value = "FORD Fiesta  Model"
$ws.Range("B8").Value = "Purchase"
$ws.Range("C8").Value = "FORD"
$ws.Range("D8").Value = "Fiesta"
$ws.Range("E8").Value = "null"
$ws.Range("F8").Value = "null"
$ws.Range("A9").Value = "FORD Focus  Model"
$ws.Range("B9").Value = "Purchase"
$ws.Range("C9").Value = "FORD"
$ws.Range("D9").Value = "Focus"
$ws.Range("E9").Value = "Kombi"
$ws.Range("F9").Value = "null"
$ws.Range("A10").Value = "FORD Kuga  Model"
$ws.Range("B10").Value = "Purchase"
$ws.Range("C10").Value = "FORD"
$ws.Range("D10").Value = "Kuga"
$ws.Range("E10").Value = "Geschlossen"
$ws.Range("F10").Value = "null"
$ws.Range("A11").Value = "Benz A class Model"
$ws.Range("B11").Value = "Purchase"
$ws.Range("C11").Value = "Mercedes-Benz"
$ws.Range("D11").Value = "A class"
$ws.Range("E11").Value = "null"
$ws.Range("F11").Value = "10.2020:0708:472"

# Column width adjustments
$ws.Columns.Item(1).ColumnWidth = 25.59
$ws.Columns.Item(6).ColumnWidth = 15.25

# Selection update
$ws.Range("D11").Select() | Out-Null
